$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")
$ws.Activate()

# Row 36: clear the "because" tagging note from column D
$ws.Range("D36").ClearContents()

# Row 37: clear the "though" tagging note from column D
$ws.Range("D37").ClearContents()

# Row 38: clear the "does not catch if" note from column D
$ws.Range("D38").ClearContents()

# Row 43: new description for how adverbs are counted
$ws.Range("C43").Value = 'Here we simply count everything that is tagged as an adverb. Biber counts all adverbs in the dictionary, everything that is longer than five letters and ends in -ly, and excludes everything that is counted as a hedge/amplifier/downtoner/placeadverbial/timeadverbial'

# Row 47: updated precision/recall note for downtoners
$ws.Range("D47").Value = 'works well'

# Row 47: new deviation note for downtoners
$ws.Range("C47").Value = 'there are some words that come to mind that could be added to Biber''s list of place adverbials: a little, a bit, a tad (HM)'

# Row 46: new precision/recall note for conjuncts
$ws.Range("D46").Value = 'doens''t work well yet. Somehow it counted "family members that are RATHER odd" even though that definitely should not be counted. It did not count "ALTOGETHER, this is nice" even though this should be counted.'

# Row 49: updated precision/recall note for amplifiers
$ws.Range("D49").Value = 'works well'

# Row 51: new precision/recall note for discourse particles
$ws.Range("D51").Value = 'works well'

# Row 68: updated precision/recall note for analytic negation
$ws.Range("D68").Value = 'works well'

# Update the view: scroll so row 42 is at the top and select D51
$win = $excel.ActiveWindow
$win.ScrollRow = 42
$win.ScrollColumn = 1
$ws.Range("D51").Select()
